{"js": "// Replace the worked-out division answers in the table with the new values.\n// Each entry maps the exact old cell text to the new cell text; all 24\n// occurrences are unique within the document, so a direct search+replace\n// for each pair is unambiguous.\nconst replacements = [\n  [\"26\u00f78=3, 2\", \"65\u00f79=7, 2\"],\n  [\"52\u00f73=17, 1\", \"59\u00f78=7, 3\"],\n  [\"12\u00f76=2, 0\", \"96\u00f77=13, 5\"],\n  [\"20\u00f79=2, 2\", \"42\u00f77=6, 0\"],\n  [\"51\u00f78=6, 3\", \"20\u00f78=2, 4\"],\n  [\"12\u00f79=1, 3\", \"69\u00f76=11, 3\"],\n  [\"24\u00f75=4, 4\", \"31\u00f73=10, 1\"],\n  [\"67\u00f79=7, 4\", \"94\u00f73=31, 1\"],\n  [\"98\u00f79=10, 8\", \"64\u00f74=16, 0\"],\n  [\"34\u00f75=6, 4\", \"74\u00f76=12, 2\"],\n  [\"56\u00f77=8, 0\", \"15\u00f78=1, 7\"],\n  [\"58\u00f78=7, 2\", \"66\u00f74=16, 2\"],\n  [\"65\u00f73=21, 2\", \"54\u00f73=18, 0\"],\n  [\"23\u00f77=3, 2\", \"89\u00f76=14, 5\"],\n  [\"39\u00f74=9, 3\", \"36\u00f75=7, 1\"],\n  [\"33\u00f75=6, 3\", \"48\u00f76=8, 0\"],\n  [\"93\u00f73=31, 0\", \"17\u00f74=4, 1\"],\n  [\"19\u00f75=3, 4\", \"90\u00f78=11, 2\"],\n  [\"59\u00f72=29, 1\", \"42\u00f79=4, 6\"],\n  [\"29\u00f76=4, 5\", \"41\u00f78=5, 1\"],\n  [\"32\u00f74=8, 0\", \"97\u00f74=24, 1\"],\n  [\"42\u00f74=10, 2\", \"89\u00f72=44, 1\"],\n  [\"84\u00f78=10, 4\", \"92\u00f74=23, 0\"],\n  [\"93\u00f74=23, 1\", \"78\u00f79=8, 6\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the worked-out division answers in the table with the new values.\n# Each entry maps the exact old cell text to the new cell text; all 24\n# occurrences are unique within the document, so Find/Replace for each\n# pair is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"26\u00f78=3, 2\";    New = \"65\u00f79=7, 2\" },\n    @{ Old = \"52\u00f73=17, 1\";   New = \"59\u00f78=7, 3\" },\n    @{ Old = \"12\u00f76=2, 0\";    New = \"96\u00f77=13, 5\" },\n    @{ Old = \"20\u00f79=2, 2\";    New = \"42\u00f77=6, 0\" },\n    @{ Old = \"51\u00f78=6, 3\";    New = \"20\u00f78=2, 4\" },\n    @{ Old = \"12\u00f79=1, 3\";    New = \"69\u00f76=11, 3\" },\n    @{ Old = \"24\u00f75=4, 4\";    New = \"31\u00f73=10, 1\" },\n    @{ Old = \"67\u00f79=7, 4\";    New = \"94\u00f73=31, 1\" },\n    @{ Old = \"98\u00f79=10, 8\";   New = \"64\u00f74=16, 0\" },\n    @{ Old = \"34\u00f75=6, 4\";    New = \"74\u00f76=12, 2\" },\n    @{ Old = \"56\u00f77=8, 0\";    New = \"15\u00f78=1, 7\" },\n    @{ Old = \"58\u00f78=7, 2\";    New = \"66\u00f74=16, 2\" },\n    @{ Old = \"65\u00f73=21, 2\";   New = \"54\u00f73=18, 0\" },\n    @{ Old = \"23\u00f77=3, 2\";    New = \"89\u00f76=14, 5\" },\n    @{ Old = \"39\u00f74=9, 3\";    New = \"36\u00f75=7, 1\" },\n    @{ Old = \"33\u00f75=6, 3\";    New = \"48\u00f76=8, 0\" },\n    @{ Old = \"93\u00f73=31, 0\";   New = \"17\u00f74=4, 1\" },\n    @{ Old = \"19\u00f75=3, 4\";    New = \"90\u00f78=11, 2\" },\n    @{ Old = \"59\u00f72=29, 1\";   New = \"42\u00f79=4, 6\" },\n    @{ Old = \"29\u00f76=4, 5\";    New = \"41\u00f78=5, 1\" },\n    @{ Old = \"32\u00f74=8, 0\";    New = \"97\u00f74=24, 1\" },\n    @{ Old = \"42\u00f74=10, 2\";   New = \"89\u00f72=44, 1\" },\n    @{ Old = \"84\u00f78=10, 4\";   New = \"92\u00f74=23, 0\" },\n    @{ Old = \"93\u00f74=23, 1\";   New = \"78\u00f79=8, 6\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2)\n}\n"}
